# Update the "想去人数" (number of people interested) figures that changed
# between the two data pulls, on both the "展览" sheet and the "全部类型"
# sheet (which contains the same exhibition rows shifted down by one row).

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 12595
$wsExhibition.Range("F5").Value = 31
$wsExhibition.Range("F8").Value = 12500
$wsExhibition.Range("F10").Value = 4913
$wsExhibition.Range("F11").Value = 4838
$wsExhibition.Range("F14").Value = 427
$wsExhibition.Range("F16").Value = 970
$wsExhibition.Range("F18").Value = 367

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F4").Value = 12595
$wsAllTypes.Range("F6").Value = 31
$wsAllTypes.Range("F9").Value = 12500
$wsAllTypes.Range("F11").Value = 4913
$wsAllTypes.Range("F12").Value = 4838
$wsAllTypes.Range("F15").Value = 427
$wsAllTypes.Range("F17").Value = 970
$wsAllTypes.Range("F19").Value = 367

$wb.Save()
